# HFF_Liquidation-8F_Air-AWB_176-50686031.xlsx
# - Add the source/lot "document" numbers to column D (pallet NO.) for rows 10-17.
#   These were previously plain numbers (e.g. 1817); the new values are the full
#   source-document numbers (e.g. 1511817) and must be stored as TEXT, not numbers,
#   matching the shared-string-backed cells produced by Excel when typed as text.
# - Update the current selection/scroll state on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $text) {
    # Preserve the cell's existing number format (e.g. "0_ ") - switching to "@"
    # only for the duration of the write is what makes Excel store the value as
    # a genuine text/string cell (t="s") while leaving the visible style (s="..")
    # on the cell unchanged once the original number format is restored.
    $originalFormat = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = $originalFormat
}

Set-TextValue $ws.Range("D10") "1511817"
Set-TextValue $ws.Range("D11") "1511818"
Set-TextValue $ws.Range("D12") "1511816"
Set-TextValue $ws.Range("D13") "1511817"
Set-TextValue $ws.Range("D14") "1511818"
Set-TextValue $ws.Range("D15") "1511819"
Set-TextValue $ws.Range("D16") "1511861"
Set-TextValue $ws.Range("D17") "1511816"

# Move the active selection/cursor to D20, matching the saved view state.
$ws.Range("D20").Select()
